$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "1.TxInquiry.TranData.TXCD`n2.TxTranCode.TranItem`n3.交易代號+交易名稱"
$ws.Range("G4").Value = "依輸入條件，輸出結果"

$ws.Range("F5").Value = "TxInquiry.TranData.REASON/查詢理由"
$ws.Range("G5").Value = "依輸入條件，輸出結果"

$ws.Range("F6").Value = "1.TxInquiry.TranData.TLRNO`n2.TxInquiry.TranData.EMPNM`n3.經辦+經辦姓名"
$ws.Range("G6").Value = "依輸入條件，輸出結果"

$ws.Range("F7").Value = "1.TxInquiry.TranData.Caldate`n2.TxInquiry.TranData.CalTime`n3.交易日期+交易時間"
$ws.Range("G7").Value = "依輸入條件，輸出結果"

[void]$ws.Range("F6").Select()
